$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Rename the sheet
# ---------------------------------------------------------------
$ws.Name = "Pilot Gantt Chart (Sprints)"

# ---------------------------------------------------------------
# 2. Column widths: A -> 90, B..M -> 15 (new columns H..M created)
#    (ColumnWidth has a constant +5/6 offset vs. the stored <col> width)
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 89.16666666666667
$ws.Range("B1:M1").EntireColumn.ColumnWidth = 14.166666666666666

# ---------------------------------------------------------------
# 3. Row 1 height 20 -> 30
# ---------------------------------------------------------------
$ws.Rows(1).RowHeight = 30

# ---------------------------------------------------------------
# 4. Header row text (A1..M1)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Activity / Task (Timeline)"
$ws.Range("B1").Value = "Sprint 1 (W1-2)"
$ws.Range("C1").Value = "Sprint 2 (W3-4)"
$ws.Range("D1").Value = "Sprint 3 (W5-6)"
$ws.Range("E1").Value = "Sprint 4 (W7-8)"
$ws.Range("F1").Value = "Sprint 5 (W9-10)"
$ws.Range("G1").Value = "Sprint 6 (W11-12)"
$ws.Range("H1").Value = "Sprint 7 (W13-14)"
$ws.Range("I1").Value = "Sprint 8 (W15-16)"
$ws.Range("J1").Value = "Sprint 9 (W17-18)"
$ws.Range("K1").Value = "Sprint 10 (W19-20)"
$ws.Range("L1").Value = "Sprint 11 (W21-22)"
$ws.Range("M1").Value = "Sprint 12 (W23-24)"

# Header cells B1:M1 get wrap text added to the existing centered style.
$ws.Range("B1:M1").WrapText = $true

# ---------------------------------------------------------------
# 5. Task row text updates (append trailing space) + re-indent style
#    (indent level 1 -> 2) applied in one pass via the Areas collection
#    so the style table only grows by a single new entry.
# ---------------------------------------------------------------
$ws.Range("A4").Value = "1. Deep Dive into Existing UAT Processes & Test Assets "
$ws.Range("A5").Value = "2. Identify & Prioritize UAT Scenarios for Automation "
$ws.Range("A6").Value = "3. Master BDD Tooling & Methodology "
$ws.Range("A8").Value = "4. Convert Selected UAT Scenarios to BDD (Gherkin) "
$ws.Range("A9").Value = "5. Develop Automated Test Scripts using Playwright "
$ws.Range("A10").Value = "6. Setup & Test Execution in DT2 Environment "
$ws.Range("A12").Value = "7. Iterate and Refine Automated UAT Suite "
$ws.Range("A13").Value = "8. Establish Automated UAT Reporting "
$ws.Range("A14").Value = "9. Document Best Practices & Create Migration Playbook "
$ws.Range("A15").Value = "10. Prepare for Knowledge Sharing & Team Onboarding "
$ws.Range("A18").Value = "1. Baseline Current Engineering Practices & CI/CD Maturity "
$ws.Range("A19").Value = "2. Develop & Communicate Pilot Engineering Practices Adoption Strategy "
$ws.Range("A20").Value = "3. Tooling Onboarding & Environment Preparation "
$ws.Range("A22").Value = "4. Drive Adoption of Unit Testing & Developer-Led Testing "
$ws.Range("A23").Value = "5. Integrate Automated Tests into CI/CD Pipelines (GitHub Actions Focus) "
$ws.Range("A24").Value = "6. Establish & Champion Mocking Practices (Mockito/MockFlow) "
$ws.Range("A26").Value = "7. Refine CI/CD Pipelines (GitHub Actions) and Test Execution Efficiency "
$ws.Range("A27").Value = "8. Develop & Document Standardized Engineering Playbooks "
$ws.Range("A28").Value = "9. Facilitate Performance Profiling Setup "
$ws.Range("A29").Value = "10. Prepare for Scaling & Knowledge Transfer "

$taskRange = $ws.Range("A4,A5,A6,A8,A9,A10,A12,A13,A14,A15,A18,A19,A20,A22,A23,A24,A26,A27,A28,A29")
foreach ($area in $taskRange.Areas) {
    $area.IndentLevel = 2
}

# ---------------------------------------------------------------
# 6. Section header rows: horizontal=left, indent=1 (previously unformatted)
# ---------------------------------------------------------------
$sectionRange = $ws.Range("A3,A7,A11,A17,A21,A25")
foreach ($area in $sectionRange.Areas) {
    $area.HorizontalAlignment = -4131
    $area.IndentLevel = 1
}

# ---------------------------------------------------------------
# 7. Remove the leftover empty cells in the section-header rows
# ---------------------------------------------------------------
$ws.Range("B3:C3").Clear()
$ws.Range("D7:E7").Clear()
$ws.Range("F11:G11").Clear()
$ws.Range("B17:C17").Clear()
$ws.Range("D21:E21").Clear()
$ws.Range("F25:G25").Clear()

# ---------------------------------------------------------------
# 8. Extend merged ranges out to column M
# ---------------------------------------------------------------
$ws.Range("A2:G2").UnMerge()
$ws.Range("A2:M2").Merge()

$ws.Range("A16:G16").UnMerge()
$ws.Range("A16:M16").Merge()

$ws.Range("A3").UnMerge()
$ws.Range("A3:M3").Merge()

$ws.Range("A7").UnMerge()
$ws.Range("A7:M7").Merge()

$ws.Range("A11").UnMerge()
$ws.Range("A11:M11").Merge()

$ws.Range("A17").UnMerge()
$ws.Range("A17:M17").Merge()

$ws.Range("A21").UnMerge()
$ws.Range("A21:M21").Merge()

$ws.Range("A25").UnMerge()
$ws.Range("A25:M25").Merge()
